$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.054.68'
$ws.Range('E2').Value = '  -0.01%  '

$ws.Range('D3').Value = '1.833.29'
$ws.Range('E3').Value = '  +0.18%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9978'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6183'
$ws.Range('E6').Value = '  -2.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07464'
$ws.Range('E8').Value = '  -1.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2928'
$ws.Range('E9').Value = '  -0.45%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.09'
$ws.Range('E10').Value = '  -0.23%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07678'
$ws.Range('E11').Value = '  -0.37%  '

$ws.Range('D12').Value = '1.833.68'
$ws.Range('E12').Value = '  +0.36%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.002'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6721'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.69'
$ws.Range('E15').Value = '  -0.61%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009176'
$ws.Range('E16').Value = '  -4.31%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.911'
$ws.Range('E17').Value = '  -2.79%  '

$ws.Range('D18').Value = '29.042.13'
$ws.Range('E18').Value = '  -0.11%  '

$ws.Range('D19').Value = '2.078.55'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '236.34'
$ws.Range('E20').Value = '  +4.27%  '

$ws.Range('E21').Value = '  +0.82%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.14%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.199'
$ws.Range('E23').Value = '  +0.76%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9995'
$ws.Range('E24').Value = '  -0.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.33'
$ws.Range('E25').Value = '  -0.70%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1423'
$ws.Range('E26').Value = '  -0.44%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.485'
$ws.Range('E27').Value = '  -0.38%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.82'
$ws.Range('E28').Value = '  -0.82%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.495'
$ws.Range('E29').Value = '  -0.72%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.153'
$ws.Range('E30').Value = '  +0.04%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.113'
$ws.Range('E31').Value = '  +1.13%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05555'
$ws.Range('E32').Value = '  +1.37%  '

$ws.Range('E33').Value = '  +0.32%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.836'
$ws.Range('E34').Value = '  -1.24%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7392'
$ws.Range('E35').Value = '  -0.75%  '

$ws.Range('E36').Value = '  -0.07%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.661'
$ws.Range('E37').Value = '  +0.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.769'
$ws.Range('E38').Value = '  +0.58%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01778'
$ws.Range('E39').Value = '  -0.36%  '

$ws.Range('D40').Value = '1.208.93'
$ws.Range('E40').Value = '  -2.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.455'
$ws.Range('E41').Value = '  -2.21%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8958'
$ws.Range('E42').Value = '  -0.84%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.02%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.82'
$ws.Range('E44').Value = '  +0.42%  '

$ws.Range('D45').Value = '1.977.74'
$ws.Range('E45').Value = '  -0.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.46'
$ws.Range('E46').Value = '  +0.64%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000120'
$ws.Range('E47').Value = '  -1.01%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5081'
$ws.Range('E48').Value = '  -0.40%  '

$ws.Range('E49').Value = '  +0.12%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.105'
$ws.Range('E50').Value = '  +1.51%  '

$ws.Range('E51').Value = '  +0.52%  '
